# Applies odds/value corrections to rows 2, 3, 4, 6, 7 of Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.6
$ws.Range("H2").Value = 2.88
$ws.Range("I2").Value = 2.35
$ws.Range("L2").Value = 3.4
$ws.Range("N2").Value = 4.75
$ws.Range("X2").Value = 15
$ws.Range("AB2").Value = 51
$ws.Range("AG2").Value = 5
$ws.Range("AH2").Value = 9
$ws.Range("AJ2").Value = 23
$ws.Range("AM2").Value = 5
$ws.Range("AV2").Value = 17

# Row 3
$ws.Range("O3").Value = 1.62
$ws.Range("P3").Value = 2.2

# Row 4
$ws.Range("M4").Value = 1.03
$ws.Range("N4").Value = 15
$ws.Range("Q4").Value = 1.7
$ws.Range("R4").Value = 2.1

# Row 6
$ws.Range("G6").Value = 3.1
$ws.Range("I6").Value = 2.4
$ws.Range("J6").Value = 3.6
$ws.Range("K6").Value = 2.1
$ws.Range("O6").Value = 1.33
$ws.Range("P6").Value = 3.4
$ws.Range("Q6").Value = 2.08
$ws.Range("R6").Value = 1.73
$ws.Range("U6").Value = 1.8
$ws.Range("V6").Value = 1.95
$ws.Range("AC6").Value = 9
$ws.Range("AO6").Value = 26
$ws.Range("AZ6").Value = 151
$ws.Range("BD6").Value = 251

# Row 7
$ws.Range("G7").Value = 2.7
$ws.Range("H7").Value = 3.1
$ws.Range("K7").Value = 1.91
$ws.Range("L7").Value = 3.6
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 7
$ws.Range("S7").Value = 1.57
$ws.Range("T7").Value = 2.25
$ws.Range("U7").Value = 2
$ws.Range("V7").Value = 1.75
$ws.Range("W7").Value = 7
$ws.Range("AA7").Value = 26
$ws.Range("AC7").Value = 6.5
$ws.Range("AE7").Value = 17
$ws.Range("AF7").Value = 67
$ws.Range("AG7").Value = 7
$ws.Range("AR7").Value = 2.25
$ws.Range("AS7").Value = 9
$ws.Range("AU7").Value = 4.75
$ws.Range("AY7").Value = 101
$ws.Range("AZ7").Value = 301
$ws.Range("BA7").Value = 301
$ws.Range("BD7").Value = 451
